# The author renamed the employee "RACINE" to "ROUX" everywhere in the
# workbook (standalone "RACINE", "Jean RACINE", and the sentences that
# mention him), then left the "Astreintes 2024" sheet active/selected at
# cell I12 and the "Astreintes 2025" sheet's last selection at F22.

$wb = $excel.ActiveWorkbook

# 1. Replace every occurrence of RACINE -> ROUX across all worksheets.
#    Doing this through Range.Replace (rather than rewriting individual
#    cell values) mirrors exactly what Excel's Find & Replace does, which
#    is what naturally reshuffles the shared-string table the same way
#    the original author's save did.
foreach ($ws in $wb.Worksheets) {
    $ws.UsedRange.Replace("RACINE", "ROUX") | Out-Null
}

# 2. Restore/update the view state: "Astreintes 2025" keeps a pending
#    selection at F22 (but is no longer the active tab), while
#    "Astreintes 2024" becomes the active sheet with its selection at I12.
$ws2025 = $wb.Worksheets.Item("Astreintes 2025")
$ws2025.Range("F22").Select() | Out-Null

$ws2024 = $wb.Worksheets.Item("Astreintes 2024")
$ws2024.Activate() | Out-Null
$ws2024.Range("I12").Select() | Out-Null
